{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per the diff):\n//  1. Remove the \"Meta description: Discover 100 Fortunes, ...\" paragraph\n//     that sits right under the \"Play 100 Fortunes Free Online Slot | See\n//     Our Review\" Heading1 at the top of the document.\n//  2. Near the very end of the document, insert a new paragraph (bold)\n//     reading \"Play 100 Fortunes Free Online Slot | See Our Review\" right\n//     before the final \"Prompt: ...\" paragraph.\n//  3. Replace the final paragraph's (italic) text \u2014 the old \"Prompt: ...\"\n//     image-generation prompt \u2014 with the meta-description copy:\n//     \"Discover 100 Fortunes, an online slot game with expanding reels and\n//     Chinese aesthetics. Play now for free and read our expert review.\"\n\nconst body = context.document.body;\n\n// --- Step 1: delete the \"Meta description\" paragraph near the top ---\nconst metaResults = body.search(\"Meta description\", { matchCase: true });\nmetaResults.load(\"items\");\nawait context.sync();\n\nif (metaResults.items.length > 0) {\n  const metaPara = metaResults.items[0].paragraphs.getFirst();\n  metaPara.delete();\n  await context.sync();\n}\n\n// --- Step 2 & 3: locate the closing \"Prompt: ...\" paragraph and replace it\n//     with the new bold heading paragraph + the updated italic paragraph ---\nconst promptResults = body.search(\"Prompt: Create a feature image\", { matchCase: true });\npromptResults.load(\"items\");\nawait context.sync();\n\nconst promptPara = promptResults.items[0].paragraphs.getFirst();\nconst wholeRange = promptPara.getRange(\"Whole\");\n\n// Flat-OPC wrapped OOXML payload \u2014 Word's insertOoxml requires the\n// \"<pkg:package>\" wrapper even for a small fragment.\nconst flatOpcXml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 100 Fortunes Free Online Slot | See Our Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover 100 Fortunes, an online slot game with expanding reels and Chinese aesthetics. Play now for free and read our expert review.</w:t></w:r></w:p><w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\nwholeRange.insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (per the diff):\n#  1. Remove the \"Meta description: Discover 100 Fortunes, ...\" paragraph\n#     that sits right under the \"Play 100 Fortunes Free Online Slot | See\n#     Our Review\" Heading1 at the top of the document.\n#  2. Near the very end of the document, insert a new paragraph (bold)\n#     reading \"Play 100 Fortunes Free Online Slot | See Our Review\" right\n#     before the final \"Prompt: ...\" paragraph.\n#  3. Replace the final paragraph's (italic) text -- the old \"Prompt: ...\"\n#     image-generation prompt -- with the meta-description copy:\n#     \"Discover 100 Fortunes, an online slot game with expanding reels and\n#     Chinese aesthetics. Play now for free and read our expert review.\"\n\n$d = $word.ActiveDocument\n\n# --- Step 1: delete the \"Meta description\" paragraph near the top ---\n$metaSearch = $d.Content\n$metaFound = $metaSearch.Find.Execute(\"Meta description\")\nif ($metaFound) {\n    $metaPara = $metaSearch.Paragraphs(1)\n    $metaPara.Range.Delete()\n}\n\n# --- Step 2 & 3: locate the closing \"Prompt: ...\" paragraph and replace it\n#     with the new bold heading paragraph + the updated italic paragraph ---\n$promptSearch = $d.Content\n$promptSearch.Find.Execute(\"Prompt: Create a feature image\") | Out-Null\n$promptPara = $promptSearch.Paragraphs(1)\n$wholeRange = $promptPara.Range\n\n# Flat-OPC wrapped OOXML payload -- Range.InsertXML requires the\n# \"<pkg:package>\" wrapper even for a small fragment.\n$flatOpcXml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 100 Fortunes Free Online Slot | See Our Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover 100 Fortunes, an online slot game with expanding reels and Chinese aesthetics. Play now for free and read our expert review.</w:t></w:r></w:p><w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$wholeRange.InsertXML($flatOpcXml)\n"}
